$p = $ppt.ActivePresentation
$newDate = "31/08/2023"

function Update-DateShapes($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master date placeholder
$m = $p.SlideMaster
Update-DateShapes $m

# Every slide layout's date placeholder
for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    Update-DateShapes $cl
}
